# Add a "Save" column (H) to the s_vals sheet, matching the existing
# header style used by the other header cells (B1:G1), and fill the
# data rows with 0 like the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cell (G1) onto H1 so the
# new header reuses the existing bold/bordered/centered header style
# instead of creating a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data column, defaulting to 0 for each existing row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
